# The Last Update 15-03-2024
# Refresh the NBA stat-leader tables (3PM, AST, PTS, REB, STL, BLK) with the
# latest values. Each sheet keeps its Rank/Nome de Jogador/Time/Valor layout;
# only the data rows (2-6) are refreshed in place.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Values like "11.4" look numeric to Excel's smart entry and would
    # otherwise be stored as a number (and lose precision). Force them to
    # stay text the same way a user would (leading apostrophe), then drop
    # the resulting "quote prefix" style so the cell's XML stays styleless,
    # matching the rest of the sheet.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 1: Arremessos de 3 Pontos (3-Pointers Made)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Set-TextValue $ws1.Range("D2") "4.9"
$ws1.Range("B5").Value = "Donovan Mitchell"
$ws1.Range("C5").Value = "CLE"
$ws1.Range("B6").Value = "Desmond Bane"
$ws1.Range("C6").Value = "MEM"

# ---------------------------------------------------------------------
# Sheet 2: Assistências (Assists)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("D2") "11.4"
Set-TextValue $ws2.Range("D3") "10.8"
Set-TextValue $ws2.Range("D4") "9.8"
Set-TextValue $ws2.Range("D5") "9.2"

# ---------------------------------------------------------------------
# Sheet 3: Pontos (Points) - Joel Embiid drops off, rows shift up
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "Luka Doncic"
$ws3.Range("C2").Value = "DAL"
Set-TextValue $ws3.Range("D2") "34.3"

$ws3.Range("B3").Value = "Shai Gilgeous-Alexander"
$ws3.Range("C3").Value = "OKC"
Set-TextValue $ws3.Range("D3") "31.1"

$ws3.Range("B4").Value = "Giannis Antetokounmpo"
$ws3.Range("C4").Value = "MIL"
Set-TextValue $ws3.Range("D4") "30.8"

$ws3.Range("B5").Value = "Kevin Durant"
$ws3.Range("C5").Value = "PHX"
Set-TextValue $ws3.Range("D5") "28.3"

$ws3.Range("B6").Value = "Donovan Mitchell"
$ws3.Range("C6").Value = "CLE"
Set-TextValue $ws3.Range("D6") "27.7"

# ---------------------------------------------------------------------
# Sheet 4: Rebotes (Rebounds)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("D2") "13.6"
Set-TextValue $ws4.Range("D3") "12.9"
Set-TextValue $ws4.Range("D4") "12.4"
Set-TextValue $ws4.Range("D5") "12.3"
Set-TextValue $ws4.Range("D6") "11.9"

# ---------------------------------------------------------------------
# Sheet 5: Roubos (Steals)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B3").Value = "De'Aaron Fox"
$ws5.Range("C3").Value = "SAC"

$ws5.Range("A4").Value = 3
$ws5.Range("B4").Value = "Matisse Thybulle"
$ws5.Range("C4").Value = "POR"
Set-TextValue $ws5.Range("D4") "1.8"

$ws5.Range("A5").Value = 3
$ws5.Range("B5").Value = "Donovan Mitchell"
$ws5.Range("C5").Value = "CLE"
Set-TextValue $ws5.Range("D5") "1.8"

$ws5.Range("B6").Value = "Kawhi Leonard"
$ws5.Range("C6").Value = "LAC"
Set-TextValue $ws5.Range("D6") "1.7"

# ---------------------------------------------------------------------
# Sheet 6: Tocos (Blocks)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("D2") "3.4"

$ws6.Range("B3").Value = "Walker Kessler"
$ws6.Range("C3").Value = "UTAH"
Set-TextValue $ws6.Range("D3") "2.6"

$ws6.Range("B4").Value = "Brook Lopez"
$ws6.Range("C4").Value = "MIL"
Set-TextValue $ws6.Range("D4") "2.6"

$ws6.Range("A5").Value = 4
Set-TextValue $ws6.Range("D5") "2.5"

Set-TextValue $ws6.Range("D6") "2.4"
